# Common: Experimental improvement of the forms
# Append new translation key/value rows to the "Import" sheet (sheet1 / rId1),
# mirroring the rows already present (Language | Key | Translation).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")

# New rows to append starting right after the last used row (856).
$startRow = 857
$endRow = 872
$lastRowStyleSource = $ws.Range("A856:C856")

for ($r = $startRow; $r -le $endRow; $r++) {
    $lastRowStyleSource.Copy()
    $ws.Range("A" + $r + ":C" + $r).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# Language column is constant.
for ($r = $startRow; $r -le $endRow; $r++) {
    $ws.Cells.Item($r, 1).Value = "cs"
}

# Key / Translation values are entered in the exact order the author typed
# them (this also governs the order new strings are added to the workbook's
# shared string table). Every row is a straightforward Key-then-Translation
# entry, except the cottonOffset.-2 / cottonOffset.-1 pair (rows 862-863),
# where the Translation cells were filled in swapped order.
$ws.Cells.Item(857, 2).Value = "lab.build.coilOffset.-2"
$ws.Cells.Item(857, 3).Value = "Nejníže"
$ws.Cells.Item(858, 2).Value = "lab.build.coilOffset.-1"
$ws.Cells.Item(858, 3).Value = "Níže"
$ws.Cells.Item(859, 2).Value = "lab.build.coilOffset.0"
$ws.Cells.Item(859, 3).Value = "Střed"
$ws.Cells.Item(860, 2).Value = "lab.build.coilOffset.1"
$ws.Cells.Item(860, 3).Value = "Výše"
$ws.Cells.Item(861, 2).Value = "lab.build.coilOffset.2"
$ws.Cells.Item(861, 3).Value = "Nejvýše"

$ws.Cells.Item(862, 2).Value = "lab.build.cottonOffset.-2"
$ws.Cells.Item(863, 3).Value = "Méně"
$ws.Cells.Item(863, 2).Value = "lab.build.cottonOffset.-1"
$ws.Cells.Item(862, 3).Value = "Nejméně"

$ws.Cells.Item(864, 2).Value = "lab.build.cottonOffset.0"
$ws.Cells.Item(864, 3).Value = "Akorát"
$ws.Cells.Item(865, 2).Value = "lab.build.cottonOffset.1"
$ws.Cells.Item(865, 3).Value = "Více"
$ws.Cells.Item(866, 2).Value = "lab.build.cottonOffset.2"
$ws.Cells.Item(866, 3).Value = "Nejvíce"
$ws.Cells.Item(867, 2).Value = "lab.build.glow.1"
$ws.Cells.Item(867, 3).Value = "Pomalé"
$ws.Cells.Item(868, 2).Value = "lab.build.glow.2"
$ws.Cells.Item(868, 3).Value = "Střední"
$ws.Cells.Item(869, 2).Value = "lab.build.glow.3"
$ws.Cells.Item(869, 3).Value = "Běžné"
$ws.Cells.Item(870, 2).Value = "lab.build.glow.4"
$ws.Cells.Item(870, 3).Value = "Rychlé"
$ws.Cells.Item(871, 2).Value = "lab.build.glow.5"
$ws.Cells.Item(871, 3).Value = "Velmi rychlé"
$ws.Cells.Item(872, 2).Value = "lab.build.preview.glow"
$ws.Cells.Item(872, 3).Value = "Rychlost žhavení"

# Update view state to reflect where the user ended up after typing.
$ws.Range("B863").Select()
